$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Piyush Chawla"
$ws.Columns.Item(1).Insert()

$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "55th"
